$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEST")

# Insert a new row above row 11 (this shifts the existing "Trauma triggers"
# rows - and everything below them - down by one). Copy row 11's formatting
# first so the newly-inserted row (and the row pushed down into position 12)
# keep the same style/height as the rest of the table.
$ws.Rows.Item(11).Copy()
$ws.Rows.Item(11).Insert()
$excel.CutCopyMode = $false

# Populate the newly inserted row 11 with the new data (a new PTSD /
# "Trauma triggers" entry about physical assault).
$ws.Range("A11").Value = "I have severe PTSD from a physical assault"
$ws.Range("B11").Value = "Would you recommend I enroll in a hands-on self-defense class?"
$ws.Range("C11").Value = "Trauma triggers"

# This row's text needs two lines to wrap at these column widths.
$ws.Rows.Item(11).RowHeight = 28

# Update the selection to match the post-edit state: whole row 11 selected.
$ws.Range("A11:XFD11").Select()
